# Apply crypto price/volume/name/link updates per the authoritative diff.
# Every text value is written via a NumberFormat="@" (Text) guard so Excel
# does not reinterpret numeric-looking strings (e.g. "63.864.39") as numbers,
# then the style is reset to "Normal" so no stray style index is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "63.864.39"  # D2
Set-TextCell 2 5 "  +0.98%  "  # E2
Set-TextCell 3 4 "3.165.82"  # D3
Set-TextCell 3 5 "  +1.70%  "  # E3
Set-TextCell 4 5 "  +0.09%  "  # E4
Set-TextCell 5 4 "588.61"  # D5
Set-TextCell 5 5 "  +0.60%  "  # E5
Set-TextCell 6 4 "146.12"  # D6
Set-TextCell 6 5 "  +0.69%  "  # E6
Set-TextCell 7 5 "  +0.06%  "  # E7
Set-TextCell 8 4 "3.158.16"  # D8
Set-TextCell 8 5 "  +1.69%  "  # E8
Set-TextCell 9 4 "0.531"  # D9
Set-TextCell 9 5 "  +0.28%  "  # E9
Set-TextCell 10 4 "0.161"  # D10
Set-TextCell 10 5 "  +6.33%  "  # E10
Set-TextCell 11 5 "  -0.83%  "  # E11
Set-TextCell 12 4 "0.461"  # D12
Set-TextCell 12 5 "  -1.53%  "  # E12
Set-TextCell 13 4 "0.0000248"  # D13
Set-TextCell 13 5 "  +0.11%  "  # E13
Set-TextCell 14 4 "37.01"  # D14
Set-TextCell 14 5 "  +3.89%  "  # E14
Set-TextCell 15 2 "WrappedliquidstakedEther2.0"  # B15
Set-TextCell 15 3 "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"  # C15
Set-TextCell 15 4 "3.687.80"  # D15
Set-TextCell 15 5 "  +1.66%  "  # E15
Set-TextCell 16 2 "TRON"  # B16
Set-TextCell 16 3 "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"  # C16
Set-TextCell 16 4 "0.122"  # D16
Set-TextCell 16 5 "  -1.27%  "  # E16
Set-TextCell 17 2 "WrappedEther"  # B17
Set-TextCell 17 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"  # C17
Set-TextCell 17 4 "3.157.11"  # D17
Set-TextCell 17 5 "  +1.39%  "  # E17
Set-TextCell 18 2 "WrappedBTC"  # B18
Set-TextCell 18 3 "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"  # C18
Set-TextCell 18 4 "63.640.24"  # D18
Set-TextCell 18 5 "  +0.75%  "  # E18
Set-TextCell 19 5 "  -0.58%  "  # E19
Set-TextCell 20 4 "466.07"  # D20
Set-TextCell 20 5 "  -0.15%  "  # E20
Set-TextCell 21 4 "14.32"  # D21
Set-TextCell 21 5 "  +1.08%  "  # E21
Set-TextCell 22 4 "0.732"  # D22
Set-TextCell 22 5 "  +0.68%  "  # E22
Set-TextCell 23 4 "7.48"  # D23
Set-TextCell 23 5 "  -0.80%  "  # E23
Set-TextCell 24 4 "13.06"  # D24
Set-TextCell 24 5 "  -2.50%  "  # E24
Set-TextCell 25 4 "81.43"  # D25
Set-TextCell 25 5 "  -0.79%  "  # E25
Set-TextCell 26 5 "  +1.34%  "  # E26
Set-TextCell 27 5 "  -0.04%  "  # E27
Set-TextCell 28 4 "9.13"  # D28
Set-TextCell 28 5 "  +6.84%  "  # E28
Set-TextCell 29 5 "  +0.51%  "  # E29
Set-TextCell 30 5 "  -0.87%  "  # E30
Set-TextCell 31 5 "  +0.05%  "  # E31
Set-TextCell 32 5 "  +2.21%  "  # E32
Set-TextCell 33 4 "27.07"  # D33
Set-TextCell 33 5 "  +0.45%  "  # E33
Set-TextCell 34 5 "  +0.27%  "  # E34
Set-TextCell 35 4 "0.0₃0863"  # D35
Set-TextCell 35 5 "  -0.69%  "  # E35
Set-TextCell 36 5 "  -0.69%  "  # E36
Set-TextCell 37 4 "3.38"  # D37
Set-TextCell 37 5 "  +1.05%  "  # E37
Set-TextCell 38 4 "2.31"  # D38
Set-TextCell 38 5 "  -4.33%  "  # E38
Set-TextCell 39 4 "6.03"  # D39
Set-TextCell 39 5 "  -1.07%  "  # E39
Set-TextCell 40 4 "50.63"  # D40
Set-TextCell 40 5 "  +0.13%  "  # E40
Set-TextCell 41 4 "443.32"  # D41
Set-TextCell 41 5 "  +1.21%  "  # E41
Set-TextCell 42 4 "8.72"  # D42
Set-TextCell 42 5 "  -0.14%  "  # E42
Set-TextCell 43 2 "Maker"  # B43
Set-TextCell 43 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"  # C43
Set-TextCell 43 4 "2.925.90"  # D43
Set-TextCell 43 5 "  +0.40%  "  # E43
Set-TextCell 44 2 "VeChain"  # B44
Set-TextCell 44 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"  # C44
Set-TextCell 44 4 "0.0372"  # D44
Set-TextCell 44 5 "  +0.37%  "  # E44
Set-TextCell 45 4 "0.274"  # D45
Set-TextCell 45 5 "  -1.71%  "  # E45
Set-TextCell 46 5 "  -1.54%  "  # E46
Set-TextCell 47 4 "36.19"  # D47
Set-TextCell 47 5 "  +3.07%  "  # E47
Set-TextCell 48 4 "125.51"  # D48
Set-TextCell 48 5 "  +2.00%  "  # E48
Set-TextCell 49 5 "  +0.03%  "  # E49
Set-TextCell 50 5 "  -0.48%  "  # E50
Set-TextCell 51 4 "24.44"  # D51
Set-TextCell 51 5 "  -0.88%  "  # E51
